# Generate Report for Handoff
# This updates the localization-status report with the latest handoff/handback
# timestamps for the "7e7efbe8-a8b7-4bc1-99b8-caba2782ea5c" entry (row 7) that
# was just (re-)handed off, as reflected across the Overview, zh-cn and de-de
# worksheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for row 7
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-03-23 16:44:42"

# zh-cn sheet: "Latest Handoff Datetime" column (E) for row 7
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-23 16:44:37"

# de-de sheet: "Latest Handoff Datetime" column (E) for row 7
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-23 16:44:42"
